$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk 0
$ws.Range("H40").Value = 2800
$ws.Range("J40").Value = 2750
$ws.Range("L40").Value = 2750
$ws.Range("N40").Value = -3100
# hunk 1
$ws.Range("H43").Value = 3000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -3138
# hunk 2
$ws.Range("H53").Value = 805.7222
$ws.Range("I53").Value = 706.9167
$ws.Range("J53").Value = 1003.3333
$ws.Range("K53").Value = 706.9167
$ws.Range("L53").Value = 1003.3333
$ws.Range("M53").Value = -69.91669999999999
$ws.Range("N53").Value = -2277.3333
# hunk 3
$ws.Range("H94").Value = 10074.75
$ws.Range("I94").Value = 10074.75
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 10074.75
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -9623.75
$ws.Range("N94").ClearContents()
# hunk 4
$ws.Range("H97").Value = 6000
$ws.Range("I97").Value = 5000
$ws.Range("J97").Value = 7000
$ws.Range("K97").Value = 15000
$ws.Range("L97").Value = 21000
$ws.Range("M97").Value = -14504
$ws.Range("N97").Value = -21992

$ws = $wb.Worksheets.Item("ARM")
# hunk 5
$ws.Range("H97").Value = 114
$ws.Range("J97").Value = 102.666664
$ws.Range("L97").Value = 102.666664
$ws.Range("N97").Value = -1094.666664
# hunk 6
$ws.Range("H132").Value = 1209.6
$ws.Range("I132").Value = 1209.6
$ws.Range("K132").Value = 3628.8
$ws.Range("M132").Value = -1098.8

$ws = $wb.Worksheets.Item("BSM")
# hunk 7
$ws.Range("H80").Value = 393.5263
$ws.Range("J80").Value = 410.58334
$ws.Range("L80").Value = 410.58334
$ws.Range("N80").Value = -2406.58334
# hunk 8
$ws.Range("H83").Value = 393.5263
$ws.Range("J83").Value = 410.58334
$ws.Range("L83").Value = 2052.9167
$ws.Range("N83").Value = -12036.9167
# hunk 9
$ws.Range("H86").Value = 1841.5385
$ws.Range("I86").Value = 1763
$ws.Range("J86").Value = 1933.1666
$ws.Range("K86").Value = 1763
$ws.Range("L86").Value = 1933.1666
$ws.Range("M86").Value = -640
$ws.Range("N86").Value = -4179.1666
# hunk 10
$ws.Range("H89").Value = 1841.5385
$ws.Range("I89").Value = 1763
$ws.Range("J89").Value = 1933.1666
$ws.Range("K89").Value = 8815
$ws.Range("L89").Value = 9665.833000000001
$ws.Range("M89").Value = -3199
$ws.Range("N89").Value = -20897.833

$ws = $wb.Worksheets.Item("CRP")
# hunk 11
$ws.Range("H5").Value = 500.5
$ws.Range("I5").Value = 165.33333
$ws.Range("J5").Value = 701.6
$ws.Range("K5").Value = 165.33333
$ws.Range("L5").Value = 701.6
$ws.Range("M5").Value = -53.33332999999999
$ws.Range("N5").Value = -925.6
# hunk 12
$ws.Range("H25").Value = 170
$ws.Range("I25").Value = 170
$ws.Range("K25").Value = 170
$ws.Range("M25").Value = 4
# hunk 13
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# hunk 14
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# hunk 15
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
# hunk 16
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
# hunk 17
$ws.Range("H107").Value = 537
$ws.Range("I107").Value = 537
$ws.Range("K107").Value = 537
$ws.Range("M107").Value = 1383

$ws = $wb.Worksheets.Item("CUL")
# hunk 18
$ws.Range("H4").Value = 125000570
$ws.Range("I4").Value = 849
$ws.Range("J4").Value = 250000290
$ws.Range("K4").Value = 2547
$ws.Range("L4").Value = 750000870
$ws.Range("M4").Value = -2435
$ws.Range("N4").Value = -750001094
# hunk 19
$ws.Range("H7").Value = 209.83333
$ws.Range("I7").Value = 197.09091
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 591.27273
$ws.Range("L7").Value = 1050
$ws.Range("M7").Value = -479.27273
$ws.Range("N7").Value = -1274
# hunk 20
$ws.Range("H52").Value = 1000
$ws.Range("J52").Value = 1000
$ws.Range("L52").Value = 3000
$ws.Range("N52").Value = -3532
# hunk 21
$ws.Range("H69").Value = 670.3333
$ws.Range("I69").Value = 670.3333
$ws.Range("K69").Value = 2010.9999
$ws.Range("M69").Value = -1199.9999
# hunk 22
$ws.Range("H72").Value = 670.3333
$ws.Range("I72").Value = 670.3333
$ws.Range("K72").Value = 6032.9997
$ws.Range("M72").Value = -1976.9997
# hunk 23
$ws.Range("H75").Value = 647.6667
$ws.Range("J75").Value = 965
$ws.Range("L75").Value = 2895
$ws.Range("N75").Value = -4891
# hunk 24
$ws.Range("H78").Value = 647.6667
$ws.Range("J78").Value = 965
$ws.Range("L78").Value = 8685
$ws.Range("N78").Value = -18669
# hunk 25
$ws.Range("H80").Value = 8193.666999999999
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 14998.25
$ws.Range("K80").Value = 8250
$ws.Range("L80").Value = 44994.75
$ws.Range("M80").Value = -7314
$ws.Range("N80").Value = -46866.75
# hunk 26
$ws.Range("H83").Value = 8193.666999999999
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 14998.25
$ws.Range("K83").Value = 24750
$ws.Range("L83").Value = 134984.25
$ws.Range("M83").Value = -20070
$ws.Range("N83").Value = -144344.25
# hunk 27
$ws.Range("H92").Value = 884.85
$ws.Range("I92").Value = 799.8333
$ws.Range("J92").Value = 921.2857
$ws.Range("K92").Value = 2399.4999
$ws.Range("L92").Value = 2763.8571
$ws.Range("M92").Value = -1151.4999
$ws.Range("N92").Value = -5259.8571
# hunk 28
$ws.Range("H109").Value = 455.6
$ws.Range("I109").Value = 332
$ws.Range("J109").Value = 950
$ws.Range("K109").Value = 996
$ws.Range("L109").Value = 2850
$ws.Range("M109").Value = 44
$ws.Range("N109").Value = -4930

$ws = $wb.Worksheets.Item("GSM")
# hunk 29
$ws.Range("H101").Value = 34995
$ws.Range("J101").Value = 34995
$ws.Range("L101").Value = 34995
$ws.Range("N101").Value = -41485
# hunk 30
$ws.Range("H122").Value = 2110.4443
$ws.Range("I122").Value = 1999.5
$ws.Range("J122").Value = 2332.3333
$ws.Range("K122").Value = 5998.5
$ws.Range("L122").Value = 6996.999899999999
$ws.Range("M122").Value = -3548.5
$ws.Range("N122").Value = -11896.9999
# hunk 31
$ws.Range("H132").Value = 2908.6316
$ws.Range("I132").Value = 2165.5715
$ws.Range("K132").Value = 6496.7145
$ws.Range("M132").Value = -3966.7145

$ws = $wb.Worksheets.Item("LTW")
# hunk 32
$ws.Range("H39").Value = 5352.6665
$ws.Range("J39").Value = 5499.5
$ws.Range("L39").Value = 5499.5
$ws.Range("N39").Value = -6419.5

$ws = $wb.Worksheets.Item("WVR")
# hunk 33
$ws.Range("H81").Value = 1250
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
# hunk 34
$ws.Range("H84").Value = 1250
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
